# The workbook gained a new run ("Run063") that was inserted into the
# "Runs" sheet at row 55, pushing the existing rows (previously 55-63,
# i.e. Run099, Run1..Run8) down by one (now 56-64). Only the first
# dozen metric columns (C:N) have data for the new run; columns O:AO
# are left blank for it, same as the XML diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Runs")

# Insert a new row at position 55; default behaviour shifts rows down.
$ws.Rows.Item(55).Insert()

# Column A carries the constant flag "0" used throughout this sheet.
$ws.Range("A55").Value = 0
$ws.Range("B55").Value = "Run063"

$ws.Range("C55").Value = "0.478 (0.442)"
$ws.Range("D55").Value = "0.634 (0.226)"
$ws.Range("E55").Value = "0.412 (0.493)"
$ws.Range("F55").Value = "0.695 (0.357)"
$ws.Range("G55").Value = "0.669 (0.152)"
$ws.Range("H55").Value = "0.714 (0.453)"
$ws.Range("I55").Value = "0.735 (0.387)"
$ws.Range("J55").Value = "0.569 (0.258)"
$ws.Range("K55").Value = "0.791 (0.407)"
$ws.Range("L55").Value = "0.604 (0.427)"
$ws.Range("M55").Value = "0.625 (0.222)"
$ws.Range("N55").Value = "0.595 (0.491)"

# Columns O55:AO55 stay blank for the new run (inserted row already
# comes back empty there, matching the self-closed <c .../> cells in
# the target XML), so nothing further to set.
